# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
# (cryptos list refresh -- Price (D) and Volume(1h) (E) updates, plus a couple
#  of coin-ranking swaps that also touch Coin (B) and Link (C)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D are stored as plain text in this sheet (some values, like
# '27.052.67', aren't valid numbers anyway). For replacement values that DO look
# like a number, force the cell to Text format first so Excel doesn't silently
# convert the literal into a numeric value.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Rows with only Price (D) and/or Volume(1h) (E) changes ---
$ws.Range("D2").Value = "27.052.67"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "1.674.02"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  +0.25%  "
Set-TextValue $ws.Range("D5") "216.22"
$ws.Range("E5").Value = "  +1.40%  "
Set-TextValue $ws.Range("D6") "0.512"
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("E9").Value = "  +1.12%  "
Set-TextValue $ws.Range("D10") "20.17"
$ws.Range("E10").Value = "  +4.78%  "
$ws.Range("E11").Value = "  +4.78%  "
$ws.Range("D12").Value = "1.909.18"
$ws.Range("E12").Value = "  +2.78%  "
$ws.Range("D13").Value = "1.675.51"
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").Value = "27.067.35"
$ws.Range("E17").Value = "  +2.10%  "
Set-TextValue $ws.Range("D18") "235.04"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("E21").Value = "  +0.18%  "
Set-TextValue $ws.Range("D22") "4.46"
$ws.Range("E22").Value = "  +3.03%  "
Set-TextValue $ws.Range("D25") "145.27"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("E29").Value = "  +0.17%  "
Set-TextValue $ws.Range("D30") "0.0498"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("E34").Value = "  +5.30%  "
$ws.Range("E35").Value = "  +5.01%  "
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  +6.75%  "
Set-TextValue $ws.Range("D39") "0.0169"
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("E40").Value = "  +3.15%  "
$ws.Range("E41").Value = "  +11.97%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("E43").Value = "  +2.92%  "
Set-TextValue $ws.Range("D44") "65.94"
$ws.Range("E44").Value = "  +4.28%  "
$ws.Range("D45").Value = "1.818.57"
$ws.Range("E45").Value = "  +2.85%  "
Set-TextValue $ws.Range("D46") "0.781"
$ws.Range("E46").Value = "  +2.43%  "
Set-TextValue $ws.Range("D47") "90.32"
$ws.Range("E47").Value = "  -0.18%  "
Set-TextValue $ws.Range("D48") "1.53"
$ws.Range("E48").Value = "  +1.20%  "

# --- Rows with Coin (B), Link (C), Price (D) and Volume(1h) (E) changes ---
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D23") "2.23"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("B24").Value = "Avalanche"
$ws.Range("C24").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D24") "9.27"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D49") "0.101"
$ws.Range("E49").Value = "  +3.98%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.0508"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "7.63"
$ws.Range("E51").Value = "  +1.07%  "
